$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.884.68'
$ws.Range("E2").Value = '  -2.92%  '
$ws.Range("D3").Value = '2.917.34'
$ws.Range("E3").Value = '  -3.68%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.79'
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.39'
$ws.Range("E6").Value = '  -4.61%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").Value = '2.916.53'
$ws.Range("E9").Value = '  -3.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.88'
$ws.Range("E10").Value = '  +5.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.144'
$ws.Range("E11").Value = '  -3.99%  '
$ws.Range("E12").Value = '  -3.82%  '
$ws.Range("E13").Value = '  -3.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.64'
$ws.Range("E14").Value = '  -5.17%  '
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").Value = '3.402.91'
$ws.Range("E16").Value = '  -3.62%  '
$ws.Range("D17").Value = '60.854.12'
$ws.Range("E17").Value = '  -2.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.76'
$ws.Range("E18").Value = '  -4.29%  '
$ws.Range("D19").Value = '2.917.88'
$ws.Range("E19").Value = '  -3.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '431.49'
$ws.Range("E20").Value = '  -4.41%  '
$ws.Range("E21").Value = '  -4.29%  '
$ws.Range("E22").Value = '  -1.73%  '
$ws.Range("E23").Value = '  -4.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.42'
$ws.Range("E24").Value = '  -3.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.83'
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("E26").Value = '  -2.98%  '
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.22'
$ws.Range("E30").Value = '  -2.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.62'
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.57'
$ws.Range("E33").Value = '  -3.35%  '
$ws.Range("E34").Value = '  -2.49%  '
$ws.Range("D35").Value = '0.0₃0875'
$ws.Range("E35").Value = '  +1.40%  '
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.65'
$ws.Range("E37").Value = '  -4.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.04'
$ws.Range("E38").Value = '  -3.68%  '
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.73'
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.01'
$ws.Range("E41").Value = '  -4.20%  '
$ws.Range("E42").Value = '  -4.49%  '
$ws.Range("E43").Value = '  -1.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.99'
$ws.Range("E44").Value = '  -2.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '377.99'
$ws.Range("E45").Value = '  -3.82%  '
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("D47").Value = '2.679.39'
$ws.Range("E47").Value = '  -1.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.54'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.73'
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("E51").Value = '  -1.66%  '
